$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Natalie's - Honey Tangerine -- quantity 2 -> 1, total cost 28.00 -> 14.00
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "1"
$ws.Range("C3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "14.00"
$ws.Range("E3").ClearFormats()

# Row 6: Natalie's - Orange Mango -- quantity 3 -> 2, total cost 39.00 -> 26.00
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "2"
$ws.Range("C6").ClearFormats()

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "26.00"
$ws.Range("E6").ClearFormats()
